$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "Sponsor Meeting" task is replaced with "Testing Old Code", with
# 0.5 hours logged on Saturday (column G) contributing to that day's total.
$ws.Range("A9").Value = "Testing Old Code"
$ws.Range("G9").Value = 0.5
$ws.Range("I9").Value = 0.5

# Row 12 ("Daily Total"): Saturday (G) and the weekly total (I) pick up the
# extra 0.5 hours worked on "Testing Old Code".
$ws.Range("G12").Value = 0.5
$ws.Range("I12").Value = 3.5

# Move the active selection to K14, matching where Caroline left the cursor
# when she saved her final timesheet for the week.
$ws.Range("K14").Select()
